$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2
$ws.Range("D2").Value = "43.803.89"
$ws.Range("E2").Value = "  +0.18%  "

# Row 3
$ws.Range("D3").Value = "2.292.05"
$ws.Range("E3").Value = "  +0.12%  "

# Row 4
$ws.Range("E4").Value = "  +0.15%  "

# Row 5
$ws.Range("D5").Value = "'114.00"
$ws.Range("E5").Value = "  +17.48%  "

# Row 6
$ws.Range("D6").Value = "'269.25"
$ws.Range("E6").Value = "  -0.17%  "

# Row 7
$ws.Range("E7").Value = "  +1.44%  "

# Row 8
$ws.Range("E8").Value = "  +0.28%  "

# Row 9
$ws.Range("E9").Value = "  +1.50%  "

# Row 10
$ws.Range("D10").Value = "'48.44"
$ws.Range("E10").Value = "  +6.98%  "

# Row 11
$ws.Range("E11").Value = "  +0.93%  "

# Row 12
$ws.Range("E12").Value = "  +13.35%  "

# Row 13
$ws.Range("E13").Value = "  +0.12%  "

# Row 14
$ws.Range("D14").Value = "'15.79"
$ws.Range("E14").Value = "  -0.57%  "

# Row 15
$ws.Range("D15").Value = "2.633.50"
$ws.Range("E15").Value = "  +0.07%  "

# Row 16
$ws.Range("D16").Value = "'0.862"
$ws.Range("E16").Value = "  +0.40%  "

# Row 17
$ws.Range("D17").Value = "2.294.05"
$ws.Range("E17").Value = "  +0.44%  "

# Row 18
$ws.Range("D18").Value = "43.658.98"
$ws.Range("E18").Value = "  -0.17%  "

# Row 19
$ws.Range("E19").Value = "  -0.94%  "

# Row 20
$ws.Range("D20").Value = "'6.88"

# Row 21
$ws.Range("D21").Value = "'72.27"
$ws.Range("E21").Value = "  +0.22%  "

# Row 22
$ws.Range("E22").Value = "  -1.25%  "

# Row 23
$ws.Range("B23").Value = "InternetComputer(DFINITY)"
$ws.Range("C23").Value = "https://coinranking.com/coin/aMNLwaUbY+internetcomputerdfinity-icp"
$ws.Range("D23").Value = "'9.96"
$ws.Range("E23").Value = "  +10.10%  "

# Row 24
$ws.Range("B24").Value = "BitcoinCash"
$ws.Range("C24").Value = "https://coinranking.com/coin/ZlZpzOJo43mIo+bitcoincash-bch"
$ws.Range("D24").Value = "'232.66"
$ws.Range("E24").Value = "  -0.04%  "

# Row 25
$ws.Range("B25").Value = "PancakeSwap"
$ws.Range("C25").Value = "https://coinranking.com/coin/ncYFcP709+pancakeswap-cake"
$ws.Range("D25").Value = "'2.97"
$ws.Range("E25").Value = "  +9.21%  "

# Row 26
$ws.Range("E26").Value = "  -0.01%  "

# Row 27
$ws.Range("D27").Value = "'11.57"
$ws.Range("E27").Value = "  +3.04%  "

# Row 28
$ws.Range("B28").Value = "LEO"
$ws.Range("C28").Value = "https://coinranking.com/coin/mqtUpyBxu8O8+leo-leo"
$ws.Range("D28").Value = "'3.90"
$ws.Range("E28").Value = "  -0.43%  "

# Row 29
$ws.Range("B29").Value = "InjectiveProtocol"
$ws.Range("C29").Value = "https://coinranking.com/coin/PkY9BmsyW+injectiveprotocol-inj"
$ws.Range("D29").Value = "'42.20"
$ws.Range("E29").Value = "  +9.52%  "

# Row 30
$ws.Range("E30").Value = "  -2.16%  "

# Row 31
$ws.Range("D31").Value = "'2.26"
$ws.Range("E31").Value = "  -1.73%  "

# Row 32
$ws.Range("D32").Value = "'175.19"
$ws.Range("E32").Value = "  -0.63%  "

# Row 33
$ws.Range("E33").Value = "  -0.93%  "

# Row 34
$ws.Range("D34").Value = "'0.0928"
$ws.Range("E34").Value = "  +3.89%  "

# Row 35
$ws.Range("E35").Value = "  +5.11%  "

# Row 36
$ws.Range("E36").Value = "  -0.16%  "

# Row 37
$ws.Range("D37").Value = "'4.76"
$ws.Range("E37").Value = "  +1.89%  "

# Row 38
$ws.Range("E38").Value = "  +3.27%  "

# Row 39
$ws.Range("D39").Value = "'0.106"
$ws.Range("E39").Value = "  -1.80%  "

# Row 40
$ws.Range("D40").Value = "'3.82"
$ws.Range("E40").Value = "  +8.96%  "

# Row 41
$ws.Range("D41").Value = "'13.89"
$ws.Range("E41").Value = "  +13.90%  "

# Row 42
$ws.Range("D42").Value = "'74.23"
$ws.Range("E42").Value = "  +15.42%  "

# Row 43
$ws.Range("D43").Value = "'2.39"
$ws.Range("E43").Value = "  +3.10%  "

# Row 44
$ws.Range("E44").Value = "  +1.66%  "

# Row 45
$ws.Range("D45").Value = "'6.34"
$ws.Range("E45").Value = "  +21.90%  "

# Row 47
$ws.Range("D47").Value = "'1.40"
$ws.Range("E47").Value = "  +3.93%  "

# Row 48
$ws.Range("E48").Value = "  +1.04%  "

# Row 49
$ws.Range("B49").Value = "TrustWalletToken"
$ws.Range("C49").Value = "https://coinranking.com/coin/Hm3OlynlC+trustwallettoken-twt"
$ws.Range("D49").Value = "'1.26"
$ws.Range("E49").Value = "  +2.79%  "

# Row 50
$ws.Range("B50").Value = "Aave"
$ws.Range("C50").Value = "https://coinranking.com/coin/ixgUfzmLR+aave-aave"
$ws.Range("D50").Value = "'102.36"
$ws.Range("E50").Value = "  +3.58%  "

# Row 51
$ws.Range("B51").Value = "Cronos"
$ws.Range("C51").Value = "https://coinranking.com/coin/65PHZTpmE55b+cronos-cro"
$ws.Range("D51").Value = "'0.0996"
$ws.Range("E51").Value = "  -2.40%  "
